# Bugfixed the naive forecaster component module
# Update the y_0_forecast (C) and y_1_forecast (E) columns:
#   - Rows 2-6: clear the forecast cells (no data available for these early rows)
#   - Rows 7-19: replace forecast values with corrected values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C2:C6 and E2:E6 (column C has no value in row 2 already, but clear anyway)
$ws.Range("C2:C6").ClearContents()
$ws.Range("E2:E6").ClearContents()

# New corrected values for rows 7-19, columns C (y_0_forecast) and E (y_1_forecast)
$values = @{
    7  = @(1.785377844167058,  2.333075171696652)
    8  = @(5.477304442308206,  4.052456259163839)
    9  = @(4.666532690711245,  3.659383764712709)
    10 = @(5.266214435142658,  4.181342739750682)
    11 = @(4.811826107786477,  4.131858242365549)
    12 = @(5.91185619417105,   4.365509285986957)
    13 = @(5.114185474093769,  5.472991335528654)
    14 = @(2.167530781895133,  2.573593955528963)
    15 = @(0.5766229317536675, 4.059584075094214)
    16 = @(2.288114387968587,  3.463553906111505)
    17 = @(-2.013802094285932, 2.374210810973465)
    18 = @(-0.5865622195987186,2.431929210693595)
    19 = @(0.7174582534189566, 2.061048937680932)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}
